$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Workbook view (xl/workbook.xml bookViews) ---
$excel.ActiveWindow.Left = 9800
$excel.ActiveWindow.Top = 460

# --- 2. New data rows 1766-1792 ---
# Column A: sector of institution (9 categories, repeated for each date block)
$sectors = @(
    "Public, 4-year or above",
    "Private nonprofit, 4-year or above",
    "Private for-profit, 4-year or above",
    "Public, 2-year",
    "Private nonprofit, 2-year",
    "Private for-profit, 2-year",
    "Public, less-than 2-year",
    "Private nonprofit, less-than 2-year",
    "Private for-profit, less-than 2-year"
)

$dates = @("sum2013", "sum2014", "sum2015")

$counts = @{
    "sum2013" = @(920, 1074, 16, 45, 4, 1, 2, $null, 4);
    "sum2014" = @(958, 1058, 19, 54, 3, 2, 1, $null, 1);
    "sum2015" = @(825, 987, 18, 89, 6, 1, 0, $null, 0)
}

$row = 1766
foreach ($d in $dates) {
    $vals = $counts[$d]
    for ($i = 0; $i -lt $sectors.Length; $i++) {
        # Set Offense (C) before Reporting Location (B) so new shared
        # strings are appended in the order "Fires - Fires" then "All".
        $ws.Cells.Item($row, 3).Value = "Fires - Fires"
        $ws.Cells.Item($row, 2).Value = "All"
        $ws.Cells.Item($row, 1).Value = $sectors[$i]
        $ws.Cells.Item($row, 4).Value = $d
        if ($null -ne $vals[$i]) {
            $ws.Cells.Item($row, 5).Value = $vals[$i]
        }
        $row++
    }
}

# Apply the same cell formatting (quote-prefixed "text" style) used by the
# rest of the sheet to columns A, B and D — this also stamps the trailing
# rows 1793-1873 (which have no values) with the same formatting only,
# matching a format-paint/paste-formats action that overran the data range.
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Range("A1766:A1873").PasteSpecial(-4122) | Out-Null
$ws.Range("B1766:B1873").PasteSpecial(-4122) | Out-Null
$ws.Range("D1766:D1873").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 4. Sheet view (top-left cell and selection) ---
$excel.ActiveWindow.ScrollRow = 1759
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E1778").Select()
